$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the sum formula to D8
$ws.Range("D8").Formula = "=B11+B10+B7+B2"

# Update the selected cell/range to D9 (matches the diff's sheetView selection)
$ws.Range("D9").Select()
